# Swap the presentation theme's colour scheme from "Integral" to the
# stock "Office Theme" palette. (In this deck the font scheme and the
# format scheme -- fills/lines/effects -- are already byte-identical
# between the "Integral" and "Office Theme" theme parts, so the only
# observable difference once the themes are swapped is the set of
# twelve theme colours.)
#
# PowerPoint's COM automation model doesn't expose a built-in RGB()
# cmdlet in this host, so provide the standard Win32 RGB() packing
# (0x00BBGGRR) ourselves.
function RGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Theme colour slots, in COM order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
$colorScheme.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1      000000
$colorScheme.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      44546A
$colorScheme.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  ED7D31
$colorScheme.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  FFC000
$colorScheme.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  4472C4
$colorScheme.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6  70AD47
$colorScheme.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink    0563C1
$colorScheme.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink 954F72
